$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 96.666664
$ws.Range("I4").Value = 96.666664
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 96.666664
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 17.333336
$ws.Range("N4").Value = $null
$ws.Range("H28").Value = 845.41174
$ws.Range("I28").Value = 737.73334
$ws.Range("J28").Value = 1653
$ws.Range("K28").Value = 737.73334
$ws.Range("L28").Value = 1653
$ws.Range("M28").Value = -252.73334
$ws.Range("N28").Value = -2623
$ws.Range("H40").Value = 1011
$ws.Range("I40").Value = 866.3333
$ws.Range("J40").Value = 1083.3334
$ws.Range("K40").Value = 866.3333
$ws.Range("L40").Value = 1083.3334
$ws.Range("M40").Value = -691.3333
$ws.Range("N40").Value = -1433.3334
$ws.Range("H53").Value = 340.92856
$ws.Range("J53").Value = 289
$ws.Range("L53").Value = 289
$ws.Range("N53").Value = -1563
$ws.Range("H64").Value = 3948.9062
$ws.Range("I64").Value = 3703.8462
$ws.Range("J64").Value = 4116.579
$ws.Range("K64").Value = 3703.8462
$ws.Range("L64").Value = 4116.579
$ws.Range("M64").Value = -3455.8462
$ws.Range("N64").Value = -4612.579
$ws.Range("H67").Value = 3948.9062
$ws.Range("I67").Value = 3703.8462
$ws.Range("J67").Value = 4116.579
$ws.Range("K67").Value = 3703.8462
$ws.Range("L67").Value = 4116.579
$ws.Range("M67").Value = -2845.8462
$ws.Range("N67").Value = -5832.579
$ws.Range("H76").Value = 7732.6333
$ws.Range("J76").Value = 5085.7144
$ws.Range("L76").Value = 5085.7144
$ws.Range("N76").Value = -5715.7144
$ws.Range("H79").Value = 7732.6333
$ws.Range("J79").Value = 5085.7144
$ws.Range("L79").Value = 5085.7144
$ws.Range("N79").Value = -7269.7144
$ws.Range("H98").Value = 1669.1578
$ws.Range("I98").Value = 1277.8823
$ws.Range("J98").Value = 4995
$ws.Range("K98").Value = 1277.8823
$ws.Range("L98").Value = 4995
$ws.Range("M98").Value = 220.1177
$ws.Range("N98").Value = -7991
$ws.Range("H122").Value = 1669.1578
$ws.Range("I122").Value = 1277.8823
$ws.Range("J122").Value = 4995
$ws.Range("K122").Value = 3833.6469
$ws.Range("L122").Value = 14985
$ws.Range("M122").Value = -1383.6469
$ws.Range("N122").Value = -19885
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2264.8147
$ws.Range("I102").Value = 1820.909
$ws.Range("J102").Value = 4218
$ws.Range("K102").Value = 1820.909
$ws.Range("L102").Value = 4218
$ws.Range("M102").Value = -198.9090000000001
$ws.Range("N102").Value = -7462
$ws.Range("H110").Value = 2167.9285
$ws.Range("I110").Value = 2261.3333
$ws.Range("J110").Value = 1999.8
$ws.Range("K110").Value = 2261.3333
$ws.Range("L110").Value = 1999.8
$ws.Range("M110").Value = -216.3332999999998
$ws.Range("N110").Value = -6089.8
$ws.Range("H132").Value = 3758
$ws.Range("I132").Value = 3660.24
$ws.Range("J132").Value = 3946
$ws.Range("K132").Value = 10980.72
$ws.Range("L132").Value = 11838
$ws.Range("M132").Value = -8450.719999999999
$ws.Range("N132").Value = -16898
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4439.3
$ws.Range("I86").Value = 2899.75
$ws.Range("J86").Value = 6748.625
$ws.Range("K86").Value = 2899.75
$ws.Range("L86").Value = 6748.625
$ws.Range("M86").Value = -1776.75
$ws.Range("N86").Value = -8994.625
$ws.Range("H89").Value = 4439.3
$ws.Range("I89").Value = 2899.75
$ws.Range("J89").Value = 6748.625
$ws.Range("K89").Value = 14498.75
$ws.Range("L89").Value = 33743.125
$ws.Range("M89").Value = -8882.75
$ws.Range("N89").Value = -44975.125
$ws.Range("H94").Value = 4988.1904
$ws.Range("I94").Value = 1069.3334
$ws.Range("J94").Value = 7927.3335
$ws.Range("K94").Value = 1069.3334
$ws.Range("L94").Value = 7927.3335
$ws.Range("M94").Value = -618.3334
$ws.Range("N94").Value = -8829.333500000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4738.4
$ws.Range("I80").Value = 5254.1665
$ws.Range("J80").Value = 2675.3333
$ws.Range("K80").Value = 5254.1665
$ws.Range("L80").Value = 2675.3333
$ws.Range("M80").Value = -4256.1665
$ws.Range("N80").Value = -4671.3333
$ws.Range("H83").Value = 4738.4
$ws.Range("I83").Value = 5254.1665
$ws.Range("J83").Value = 2675.3333
$ws.Range("K83").Value = 26270.8325
$ws.Range("L83").Value = 13376.6665
$ws.Range("M83").Value = -21278.8325
$ws.Range("N83").Value = -23360.6665
$ws.Range("H102").Value = 3791.3928
$ws.Range("I102").Value = 4115.0435
$ws.Range("J102").Value = 2302.6
$ws.Range("K102").Value = 4115.0435
$ws.Range("L102").Value = 2302.6
$ws.Range("M102").Value = -2493.0435
$ws.Range("N102").Value = -5546.6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 6728.5713
$ws.Range("I26").Value = 5945
$ws.Range("J26").Value = 7042
$ws.Range("K26").Value = 5945
$ws.Range("L26").Value = 7042
$ws.Range("M26").Value = -5650
$ws.Range("N26").Value = -7632
$ws.Range("H29").Value = 6844
$ws.Range("I29").Value = 5833.3335
$ws.Range("J29").Value = 7602
$ws.Range("K29").Value = 5833.3335
$ws.Range("L29").Value = 7602
$ws.Range("M29").Value = -5538.3335
$ws.Range("N29").Value = -8192
$ws.Range("H61").Value = 1896.5385
$ws.Range("I61").Value = 1959.8182
$ws.Range("J61").Value = 1548.5
$ws.Range("K61").Value = 1959.8182
$ws.Range("L61").Value = 1548.5
$ws.Range("M61").Value = -1757.8182
$ws.Range("N61").Value = -1952.5
$ws.Range("H113").Value = 1896.5385
$ws.Range("I113").Value = 1959.8182
$ws.Range("J113").Value = 1548.5
$ws.Range("K113").Value = 1959.8182
$ws.Range("L113").Value = 1548.5
$ws.Range("M113").Value = 210.1818000000001
$ws.Range("N113").Value = -5888.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 12666.667
$ws.Range("I32").Value = 8000
$ws.Range("K32").Value = 8000
$ws.Range("M32").Value = -7683
$ws.Range("H34").Value = 7882.5713
$ws.Range("I34").Value = 5035.6
$ws.Range("K34").Value = 5035.6
$ws.Range("M34").Value = -4832.6
$ws.Range("H122").Value = 43307.082
$ws.Range("I122").Value = 57095
$ws.Range("J122").Value = 1943.3334
$ws.Range("K122").Value = 171285
$ws.Range("L122").Value = 5830.0002
$ws.Range("M122").Value = -168835
$ws.Range("N122").Value = -10730.0002
